$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that currently sits right after the
#    "The overall goal of this scenario..." paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the "Break the problem apart" paragraph and insert two new
#    sub-bullet paragraphs right after it.
# ------------------------------------------------------------------
$breakApartIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd([char]13, [char]7) -eq "Break the problem apart") {
        $breakApartIndex = $i
        break
    }
}

$breakApartPara = $d.Paragraphs.Item($breakApartIndex)
$breakApartPara.Range.InsertParagraphAfter()

# --- First new paragraph: "The constraints the man faces are ..." ---
$para1 = $d.Paragraphs.Item($breakApartIndex + 1)
$para1Start = $para1.Range.Start
$firstRunText = "The constraints the man faces are the "
$secondRunText = "inability to leave the parrot with the bag of seed, the likelihood of the cat eating the parrot, and the ability to only transport one item per trip."

$para1.Range.InsertBefore($firstRunText)
$splitPoint = $para1Start + $firstRunText.Length
$para1.Range.InsertAfter($secondRunText)

# Promote this paragraph to the second list level (ilvl=1), matching
# the other sub-bullets under "Define the problem".
$para1.Range.ListFormat.ListLevelNumber = 2

# The paragraph-mark formatting should stay bold (inherited from
# "Break the problem apart"); only the actual run text should not be
# bold, so operate on a tight Range that excludes the mark.
$para1End = $para1.Range.End
$run1Range = $d.Range($para1Start, $splitPoint)
$run1Range.Font.Bold = $false
$run2Range = $d.Range($splitPoint, $para1End)
$run2Range.Font.Bold = $false

# --- Second new paragraph: "The sub-goals in this scenario ..." ---
$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item($breakApartIndex + 2)
$para2Text = "The sub-goals in this scenario are to safely transport each item to the opposite riverbank, and to not leave the wrong items together while transporting the third in his boat."
$para2.Range.InsertBefore($para2Text)
$para2.Range.ListFormat.ListLevelNumber = 2

$para2Start = $para2.Range.Start
$para2End = $para2.Range.End
$run3Range = $d.Range($para2Start, $para2End)
$run3Range.Font.Bold = $false

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark at the very end of this second
#    new paragraph (collapsed, zero-length, right before the
#    paragraph mark) -- matching the original markup's placement.
#
#    A collapsed Range sitting exactly at (paragraph-end - 1) is
#    mishandled by this host's Bookmarks.Add (it silently snaps to
#    the start of the paragraph instead), so work around it: insert a
#    throw-away character, bookmark the range spanning it, then
#    delete the character again. The bookmark collapses back down to
#    the correct position and survives the deletion.
# ------------------------------------------------------------------
$para2 = $d.Paragraphs.Item($breakApartIndex + 2)
$pos = $para2.Range.End - 1
$tempRange = $d.Range($pos, $pos)
$tempRange.InsertAfter("X")
$bmRange = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$delRange = $d.Range($pos, $pos + 1)
$delRange.Text = ""
